$wb = $excel.ActiveWorkbook

# ---- Workbook-level: move the active / selected tab from Mamatha to Ram ----
$ram = $wb.Worksheets.Item("Ram")
$ram.Activate() | Out-Null

# ---- Fill in the new tracker rows on the "Ram" sheet ----

# Style donor cells already present in the workbook that carry the exact
# formatting the new rows need (reusing them keeps the style table free of
# redundant duplicate entries, same as Excel would do).
$dateStyleSrc   = $wb.Worksheets.Item("Monisha").Range("B2")   # centered date format
$taskStyleSrc   = $wb.Worksheets.Item("Monisha").Range("C2")   # centered text, border
$noFillBoxSrc   = $ram.Range("C2")                             # plain boxed cell (no fill / numFmt 0)
$plainBoxSrc    = $ram.Range("A2")                             # plain boxed "No" cell style

$rows = @(
    @{ Row=2; No=1; Date=44067; App='Vastu'; Task='isuue fixing and add DB value in japanese language '; Status='WIP' },
    @{ Row=3; No=2; Date=44068; App='Vastu'; Task='isuue fixing and add DB in japanese language '; Status='completed' },
    @{ Row=4; No=3; Date=44069; App='Vastu'; Task='isuue fixing vastu tips and homevastu'; Status='completed' },
    @{ Row=5; No=4; Date=44070; App='Vastu'; Task='isuue fixing in vastu'; Status='WIP' },
    @{ Row=6; No=5; Date=44071; App='Vastu'; Task='isuue fixing in vastu'; Status='WIP' }
)

foreach ($r in $rows) {
    $row = $r.Row

    $plainBoxSrc.Copy()
    $ram.Range("A$row").PasteSpecial(-4122)
    $ram.Range("A$row").Value = $r.No

    $dateStyleSrc.Copy()
    $ram.Range("B$row").PasteSpecial(-4122)
    $ram.Range("B$row").Value = $r.Date

    $taskStyleSrc.Copy()
    $ram.Range("C$row").PasteSpecial(-4122)
    $ram.Range("C$row").Value = $r.App

    # Wrapped, vertically centered comment cell (border all sides)
    $ram.Range("D$row").Borders.LineStyle = 1
    $ram.Range("D$row").WrapText = $true
    $ram.Range("D$row").VerticalAlignment = -4108
    $ram.Range("D$row").Value = $r.Task

    # Centered cell with no top border (continuation look)
    $noFillBoxSrc.Copy()
    $ram.Range("E$row").PasteSpecial(-4122)
    $ram.Range("E$row").Borders.Item(8).LineStyle = -4142
    $ram.Range("E$row").HorizontalAlignment = -4108
    $ram.Range("E$row").VerticalAlignment = -4108

    $noFillBoxSrc.Copy()
    $ram.Range("F$row").PasteSpecial(-4122)
    $ram.Range("F$row").Borders.Item(8).LineStyle = -4142
    $ram.Range("F$row").HorizontalAlignment = -4108
    $ram.Range("F$row").VerticalAlignment = -4108
    $ram.Range("F$row").Value = $r.Status
}

# Rows 5 & 6 also gain a boxed (empty) G cell, matching the other data rows
$ram.Range("G5").Borders.LineStyle = 1
$ram.Range("G6").Borders.LineStyle = 1

# ---- Final selection: Ram!F12 becomes the active cell/tab ----
$ram.Range("F12").Select() | Out-Null
